# Update the contact list: re-capitalize names, shift the email list down,
# drop the "goutham" row's first name, and add Gurijala/Vamsi as new rows
# with a fresh hyperlink for the new email address.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Sandeep Konduri (last name re-capitalized)
$ws.Range("B2").Value = "Konduri"

# Row 3 - Santhosh Konduri (first + last name re-capitalized)
$ws.Range("A3").Value = "Santhosh"
$ws.Range("B3").Value = "Konduri"

# Row 4 - first name cleared, last name + email replaced
$ws.Range("A4").Value = ""
$ws.Range("B4").Value = "Gurijala"
$ws.Range("C4").Value = "vamsi.gch@gmail.com"

# Row 5 - new first name, new email with its own hyperlink
$ws.Range("A5").Value = "Vamsi"
$ws.Range("C5").Value = "gurijala2018@gmail.com"
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:gurijala2018@gmail.com")
$ws.Range("C5").Style = "Hyperlink"
